$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '27.659.60'
Set-TextValue $ws.Range("E2") '  +1.75%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.766.75'
Set-TextValue $ws.Range("E3") '  -1.12%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.002'
Set-TextValue $ws.Range("E4") '  -0.33%  '

# Row 5
Set-TextValue $ws.Range("D5") '335.82'
Set-TextValue $ws.Range("E5") '  -0.17%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.9991'
Set-TextValue $ws.Range("E6") '  -0.33%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.3840'
Set-TextValue $ws.Range("E7") '  +0.06%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.3420'
Set-TextValue $ws.Range("E8") '  -0.36%  '

# Row 9
Set-TextValue $ws.Range("D9") '47.04'
Set-TextValue $ws.Range("E9") '  -2.80%  '

# Row 10
Set-TextValue $ws.Range("D10") '1.138'
Set-TextValue $ws.Range("E10") '  -5.06%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.07406'
Set-TextValue $ws.Range("E11") '  -1.47%  '

# Row 12
Set-TextValue $ws.Range("D12") '1.000'
Set-TextValue $ws.Range("E12") '  -0.31%  '

# Row 13
Set-TextValue $ws.Range("D13") '22.31'
Set-TextValue $ws.Range("E13") '  +1.97%  '

# Row 14
Set-TextValue $ws.Range("D14") '6.351'
Set-TextValue $ws.Range("E14") '  -1.87%  '

# Row 15
Set-TextValue $ws.Range("D15") '1.767.04'
Set-TextValue $ws.Range("E15") '  -1.30%  '

# Row 16
Set-TextValue $ws.Range("D16") '7.082'
Set-TextValue $ws.Range("E16") '  -0.20%  '

# Row 17
Set-TextValue $ws.Range("D17") '0.00001073'
Set-TextValue $ws.Range("E17") '  -2.08%  '

# Row 18
Set-TextValue $ws.Range("D18") '0.06654'
Set-TextValue $ws.Range("E18") '  -0.37%  '

# Row 19
Set-TextValue $ws.Range("E19") '  -2.23%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.9992'
Set-TextValue $ws.Range("E20") '  -0.31%  '

# Row 21
Set-TextValue $ws.Range("D21") '17.32'
Set-TextValue $ws.Range("E21") '  -0.45%  '

# Row 22
Set-TextValue $ws.Range("D22") '6.400'
Set-TextValue $ws.Range("E22") '  -3.56%  '

# Row 23
Set-TextValue $ws.Range("D23") '27.649.59'
Set-TextValue $ws.Range("E23") '  +1.69%  '

# Row 24
Set-TextValue $ws.Range("D24") '12.03'
Set-TextValue $ws.Range("E24") '  -2.67%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.385'
Set-TextValue $ws.Range("E25") '  +0.02%  '

# Row 26
Set-TextValue $ws.Range("D26") '1.426'
Set-TextValue $ws.Range("E26") '  -3.17%  '

# Row 27
Set-TextValue $ws.Range("D27") '20.63'
Set-TextValue $ws.Range("E27") '  -3.30%  '

# Row 28
Set-TextValue $ws.Range("D28") '2.412'
Set-TextValue $ws.Range("E28") '  -5.23%  '

# Row 29
Set-TextValue $ws.Range("D29") '152.09'
Set-TextValue $ws.Range("E29") '  -1.50%  '

# Row 30
Set-TextValue $ws.Range("D30") '134.29'
Set-TextValue $ws.Range("E30") '  -0.31%  '

# Row 31
Set-TextValue $ws.Range("D31") '1.966.42'
Set-TextValue $ws.Range("E31") '  -1.27%  '

# Row 32
Set-TextValue $ws.Range("D32") '6.118'
Set-TextValue $ws.Range("E32") '  +0.49%  '

# Row 33
Set-TextValue $ws.Range("D33") '3.962'
Set-TextValue $ws.Range("E33") '  -1.38%  '

# Row 34
Set-TextValue $ws.Range("D34") '0.08798'
Set-TextValue $ws.Range("E34") '  +0.92%  '

# Row 35
Set-TextValue $ws.Range("D35") '12.73'
Set-TextValue $ws.Range("E35") '  -4.24%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.02415'
Set-TextValue $ws.Range("E36") '  +3.12%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.6786'
Set-TextValue $ws.Range("E37") '  -1.93%  '

# Row 38
Set-TextValue $ws.Range("D38") '5.318'
Set-TextValue $ws.Range("E38") '  -2.26%  '

# Row 39
Set-TextValue $ws.Range("D39") '0.06304'
Set-TextValue $ws.Range("E39") '  -0.79%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.2179'
Set-TextValue $ws.Range("E40") '  -1.19%  '

# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D41") '1.250'
Set-TextValue $ws.Range("E41") '  +0.81%  '

# Row 42
$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D42") '1.510'
Set-TextValue $ws.Range("E42") '  -8.65%  '

# Row 43
Set-TextValue $ws.Range("D43") '8.262'
Set-TextValue $ws.Range("E43") '  -6.05%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D44") '14.19'
Set-TextValue $ws.Range("E44") '  -1.44%  '

# Row 45
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D45") '0.9989'
Set-TextValue $ws.Range("E45") '  -0.31%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.6246'
Set-TextValue $ws.Range("E46") '  -3.88%  '

# Row 47
Set-TextValue $ws.Range("D47") '3.847'
Set-TextValue $ws.Range("E47") '  -0.16%  '

# Row 48
Set-TextValue $ws.Range("D48") '131.66'
Set-TextValue $ws.Range("E48") '  +1.04%  '

# Row 49
Set-TextValue $ws.Range("D49") '2.071'
Set-TextValue $ws.Range("E49") '  -3.62%  '

# Row 50
Set-TextValue $ws.Range("D50") '0.07398'
Set-TextValue $ws.Range("E50") '  +3.67%  '

# Row 51
$ws.Range("B51").Value = 'Tezos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz'
Set-TextValue $ws.Range("D51") '1.233'
Set-TextValue $ws.Range("E51") '  +2.03%  '

Write-Output "Applied cryptos update"